# Adds snowboard variant rows 4-11 (snowboard102 .. snowboard109) to the
# ProductCards sheet, extending the used range from A1:R3 to A1:R11.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row 4
$ws.Range("A4").Value = "variant"
$ws.Range("B4").Value = "snowboard102"
$ws.Range("C4").Value = "snowboardTypeKey"
$ws.Range("D4").Value = "product-type"
$ws.Range("E4").Value = "CAPiTA Aeronaut Snowboard 2025"
$ws.Range("F4").Value = "snowboardSlug102"
$ws.Range("G4").Value = "SNW-2-01"
$ws.Range("H4").Value = "'153"
$ws.Range("I4").Value = "'1200"
$ws.Range("J4").Value = "'292"
$ws.Range("K4").Value = "'247"
$ws.Range("L4").Value = "'287"
$ws.Range("M4").Value = "8.3 ​/ 7.3"
$ws.Range("N4").Value = "'"
$ws.Range("O4").Value = "'21"
$ws.Range("P4").Value = "Regular"
$ws.Range("Q4").Value = "Intermediate-AdvancedKey;CAPiTAKey;adult-maleKey;All-MountainKey"
$ws.Range("R4").Value = "115 - 175"

# row 5
$ws.Range("A5").Value = "variant"
$ws.Range("B5").Value = "snowboard103"
$ws.Range("C5").Value = "snowboardTypeKey"
$ws.Range("D5").Value = "product-type"
$ws.Range("E5").Value = "Season Nexus Snowboard"
$ws.Range("F5").Value = "snowboardSlug103"
$ws.Range("G5").Value = "SNW-3-01"
$ws.Range("H5").Value = "'143"
$ws.Range("I5").Value = "'1030"
$ws.Range("J5").Value = "'278"
$ws.Range("K5").Value = "'242"
$ws.Range("L5").Value = "'272"
$ws.Range("M5").Value = "'7"
$ws.Range("N5").Value = "'20"
$ws.Range("O5").Value = "18 - 23"
$ws.Range("P5").Value = "Regular"
$ws.Range("Q5").Value = "Intermediate-AdvancedKey;SeasonKey;adult-anyKey;All-MountainKey;FreerideKey"
$ws.Range("R5").Value = "'"

# row 6
$ws.Range("A6").Value = "variant"
$ws.Range("B6").Value = "snowboard104"
$ws.Range("C6").Value = "snowboardTypeKey"
$ws.Range("D6").Value = "product-type"
$ws.Range("E6").Value = "Lib Tech T.Rice Orca Snowboard 2024"
$ws.Range("F6").Value = "snowboardSlug104"
$ws.Range("G6").Value = "SNW-4-01"
$ws.Range("H6").Value = "'138"
$ws.Range("I6").Value = "'1000"
$ws.Range("J6").Value = "'288"
$ws.Range("K6").Value = "'247"
$ws.Range("L6").Value = "'280"
$ws.Range("M6").Value = "'6.8"
$ws.Range("N6").Value = "'2.0"
$ws.Range("O6").Value = "17.5-22.5"
$ws.Range("P6").Value = "Regular"
$ws.Range("Q6").Value = "Intermediate-AdvancedKey;LibTechKey;adult-anyKey;All-MountainKey;PowderKey"
$ws.Range("R6").Value = "'"

# row 7
$ws.Range("A7").Value = "variant"
$ws.Range("B7").Value = "snowboard105"
$ws.Range("C7").Value = "snowboardTypeKey"
$ws.Range("D7").Value = "product-type"
$ws.Range("E7").Value = "Rossignol Myth Snowboard - Women's 2024"
$ws.Range("F7").Value = "snowboardSlug105"
$ws.Range("G7").Value = "SNW-5-01"
$ws.Range("H7").Value = "'139"
$ws.Range("I7").Value = "'1060"
$ws.Range("J7").Value = "'273"
$ws.Range("K7").Value = "'234"
$ws.Range("L7").Value = "'273"
$ws.Range("M7").Value = "'6.7"
$ws.Range("N7").Value = "'0"
$ws.Range("O7").Value = "17.3 - 20.5"
$ws.Range("P7").Value = "'"
$ws.Range("Q7").Value = "Beginner-IntermediateKey;RossignolKey;adult-femaleKey;FreestyleKey"
$ws.Range("R7").Value = "65 - 110"

# row 8
$ws.Range("A8").Value = "variant"
$ws.Range("B8").Value = "snowboard106"
$ws.Range("C8").Value = "snowboardTypeKey"
$ws.Range("D8").Value = "product-type"
$ws.Range("E8").Value = "Ride Warpig Snowboard 2024"
$ws.Range("F8").Value = "snowboardSlug106"
$ws.Range("G8").Value = "SNW-6-01"
$ws.Range("H8").Value = "'142"
$ws.Range("I8").Value = "'"
$ws.Range("J8").Value = "'301"
$ws.Range("K8").Value = "'250"
$ws.Range("L8").Value = "'301"
$ws.Range("M8").Value = "4.6​/5.6"
$ws.Range("N8").Value = "'"
$ws.Range("O8").Value = "'"
$ws.Range("P8").Value = "Regular"
$ws.Range("Q8").Value = "Intermediate-AdvancedKey;RideKey;adult-anyKey;All-MountainKey"
$ws.Range("R8").Value = "60-160"

# row 9
$ws.Range("A9").Value = "variant"
$ws.Range("B9").Value = "snowboard107"
$ws.Range("C9").Value = "snowboardTypeKey"
$ws.Range("D9").Value = "product-type"
$ws.Range("E9").Value = "Lib Tech Cold Brew C2 Snowboard 2024"
$ws.Range("F9").Value = "snowboardSlug107"
$ws.Range("G9").Value = "SNW-7-01"
$ws.Range("H9").Value = "'149"
$ws.Range("I9").Value = "'1110"
$ws.Range("J9").Value = "'295"
$ws.Range("K9").Value = "'256"
$ws.Range("L9").Value = "'292"
$ws.Range("M9").Value = "'8.1"
$ws.Range("N9").Value = "'1.5"
$ws.Range("O9").Value = "'5.5"
$ws.Range("P9").Value = "Regular"
$ws.Range("Q9").Value = "Intermediate-AdvancedKey;LibTechKey;adult-maleKey;All-MountainKey"
$ws.Range("R9").Value = "'"

# row 10
$ws.Range("A10").Value = "variant"
$ws.Range("B10").Value = "snowboard108"
$ws.Range("C10").Value = "snowboardTypeKey"
$ws.Range("D10").Value = "product-type"
$ws.Range("E10").Value = "K2 Dreamsicle Snowboard - Women's 2024"
$ws.Range("F10").Value = "snowboardSlug108"
$ws.Range("G10").Value = "SNW-8-01"
$ws.Range("H10").Value = "'138"
$ws.Range("I10").Value = "'1080"
$ws.Range("J10").Value = "'275"
$ws.Range("K10").Value = "'237"
$ws.Range("L10").Value = "'275"
$ws.Range("M10").Value = "'7.1"
$ws.Range("N10").Value = "'0.75"
$ws.Range("O10").Value = "'18"
$ws.Range("P10").Value = "Regular"
$ws.Range("Q10").Value = "Beginner-IntermediateKey;K2Key;adult-femaleKey;All-MountainKey"
$ws.Range("R10").Value = "90-160​+"

# row 11
$ws.Range("A11").Value = "variant"
$ws.Range("B11").Value = "snowboard109"
$ws.Range("C11").Value = "snowboardTypeKey"
$ws.Range("D11").Value = "product-type"
$ws.Range("E11").Value = "CAPiTA Mega Mercury Snowboard 2024"
$ws.Range("F11").Value = "snowboardSlug109"
$ws.Range("G11").Value = "SNW-9-01"
$ws.Range("H11").Value = "'153"
$ws.Range("I11").Value = "'1183"
$ws.Range("J11").Value = "'296"
$ws.Range("K11").Value = "'263"
$ws.Range("L11").Value = "'296"
$ws.Range("M11").Value = "7.5 ​/ 1.5 ​/ 7.5"
$ws.Range("N11").Value = "'0.5"
$ws.Range("O11").Value = "'529"
$ws.Range("P11").Value = "Regular"
$ws.Range("Q11").Value = "Advanced-ExpertKey;CAPiTAKey;adult-maleKey;All-MountainKey;FreerideKey"
$ws.Range("R11").Value = "100-160"

